# Adds test cases TC002-TC008 to the Adactin Login test sheet (Sheet1),
# mirroring the TC001 block already present in rows 7-12.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# TC002 - Invalid username, valid password
$ws.Range("A13").Value = "Adactin_Login_TS001_TC002"
$ws.Range("B13").Value = "To validdate login using Invalid credentials (valid username and Invalid password)"
$ws.Range("B13").WrapText = $true
$ws.Range("C13").Value = "Adactin_Login_TS001"
$ws.Range("D13").Value = "user should have registered"
$ws.Range("D13").WrapText = $true
$ws.Range("E13").Value = "Step 1"
$ws.Range("F13").Value = "start  browser"
$ws.Range("G13").Value = "chrome"
$ws.Range("H13").Value = "browser should get started"
$ws.Rows.Item(13).RowHeight = 117

$ws.Range("E14").Value = "Step 2"
$ws.Range("F14").Value = "launch app"
$ws.Range("G14").Value = "https://adactinhotelapp.com/"
$ws.Range("G14").Style = "Hyperlink"
$ws.Range("H14").Value = "app should get launched and land in home page"

$ws.Range("E15").Value = "Step 3"
$ws.Range("F15").Value = "enter username in username textbox"
$ws.Range("G15").Value = "reyaz0806"
$ws.Range("H15").Value = "textbox should accept the text and display the same"

$ws.Range("E16").Value = "Step 4"
$ws.Range("F16").Value = "enter password in password textbox"
$ws.Range("G16").Value = "reyaz456"
$ws.Range("H16").Value = "textbox should accept the text and display  as masked"

$ws.Range("E17").Value = "Step 5"
$ws.Range("F17").Value = "click login button"
$ws.Range("H17").Value = "Error shoud display below message : `n`"`t`nInvalid Login details or Your Password might have expired. Click here to reset your password`""
$ws.Range("H17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 93.6

# TC003 - Valid username, invalid password
$ws.Range("A19").Value = "Adactin_Login_TS001_TC003"
$ws.Range("B19").Value = "To validdate login using Invalid credentials (Invalid username and valid password)"
$ws.Range("B19").WrapText = $true
$ws.Range("C19").Value = "Adactin_Login_TS001"
$ws.Range("D19").Value = "user should have registered"
$ws.Range("D19").WrapText = $true
$ws.Range("E19").Value = "Step 1"
$ws.Range("F19").Value = "start  browser"
$ws.Range("G19").Value = "chrome"
$ws.Range("H19").Value = "browser should get started"
$ws.Rows.Item(19).RowHeight = 117

$ws.Range("E20").Value = "Step 2"
$ws.Range("F20").Value = "launch app"
$ws.Range("G20").Value = "https://adactinhotelapp.com/"
$ws.Range("G20").Style = "Hyperlink"
$ws.Range("H20").Value = "app should get launched and land in home page"

$ws.Range("E21").Value = "Step 3"
$ws.Range("F21").Value = "enter username in username textbox"
$ws.Range("G21").Value = "reyaz0602"
$ws.Range("H21").Value = "textbox should accept the text and display the same"

$ws.Range("E22").Value = "Step 4"
$ws.Range("F22").Value = "enter password in password textbox"
$ws.Range("G22").Value = "reyaz123"
$ws.Range("H22").Value = "textbox should accept the text and display  as masked"

$ws.Range("E23").Value = "Step 5"
$ws.Range("F23").Value = "click login button"
$ws.Range("H23").Value = "Error shoud display below message : `n`"`t`nInvalid Login details or Your Password might have expired. Click here to reset your password`""
$ws.Range("H23").WrapText = $true
$ws.Rows.Item(23).RowHeight = 93.6

# TC004 - Invalid username, invalid password
$ws.Range("A25").Value = "Adactin_Login_TS001_TC004"
$ws.Range("B25").Value = "To validdate login using Invalid credentials (Invalid username and Invalid password)"
$ws.Range("B25").WrapText = $true
$ws.Range("C25").Value = "Adactin_Login_TS001"
$ws.Range("D25").Value = "user should have registered"
$ws.Range("D25").WrapText = $true
$ws.Range("E25").Value = "Step 1"
$ws.Range("F25").Value = "start  browser"
$ws.Range("G25").Value = "chrome"
$ws.Range("H25").Value = "browser should get started"
$ws.Rows.Item(25).RowHeight = 117

$ws.Range("E26").Value = "Step 2"
$ws.Range("F26").Value = "launch app"
$ws.Range("G26").Value = "https://adactinhotelapp.com/"
$ws.Range("G26").Style = "Hyperlink"
$ws.Range("H26").Value = "app should get launched and land in home page"

$ws.Range("E27").Value = "Step 3"
$ws.Range("F27").Value = "enter username in username textbox"
$ws.Range("G27").Value = "reyaz0602"
$ws.Range("H27").Value = "textbox should accept the text and display the same"

$ws.Range("E28").Value = "Step 4"
$ws.Range("F28").Value = "enter password in password textbox"
$ws.Range("G28").Value = "reyaz456"
$ws.Range("H28").Value = "textbox should accept the text and display  as masked"

$ws.Range("E29").Value = "Step 5"
$ws.Range("F29").Value = "click login button"
$ws.Range("H29").Value = "Error shoud display below message : `n`"`t`nInvalid Login details or Your Password might have expired. Click here to reset your password`""
$ws.Range("H29").WrapText = $true
$ws.Rows.Item(29).RowHeight = 93.6

# TC005 - Blank username and password
$ws.Range("A31").Value = "Adactin_Login_TS001_TC005"
$ws.Range("B31").Value = "To validdate login with both username and password  blank fields"
$ws.Range("B31").WrapText = $true
$ws.Range("C31").Value = "Adactin_Login_TS001"
$ws.Range("D31").Value = "NA"
$ws.Range("E31").Value = "Step 1"
$ws.Range("F31").Value = "start  browser"
$ws.Range("G31").Value = "chrome"
$ws.Rows.Item(31).RowHeight = 93.6

$ws.Range("E32").Value = "Step 2"
$ws.Range("F32").Value = "launch app"
$ws.Range("G32").Value = "https://adactinhotelapp.com/"
$ws.Range("G32").Style = "Hyperlink"

$ws.Range("E33").Value = "Step 3"
$ws.Range("F33").Value = "click login button"
$ws.Range("H33").Value = "Error should display below text `n`"Enter Username`""
$ws.Range("H33").WrapText = $true
$ws.Rows.Item(33).RowHeight = 46.8

# TC006 - Blank username, password filled
$ws.Range("A35").Value = "Adactin_Login_TS001_TC006"
$ws.Range("B35").Value = "To validdate login with username and password  blank field"
$ws.Range("B35").WrapText = $true
$ws.Range("C35").Value = "Adactin_Login_TS001"
$ws.Range("D35").Value = "user should have registered"
$ws.Range("D35").WrapText = $true
$ws.Range("E35").Value = "Step 1"
$ws.Range("F35").Value = "start  browser"
$ws.Range("G35").Value = "chrome"
$ws.Rows.Item(35).RowHeight = 70.2

$ws.Range("E36").Value = "Step 2"
$ws.Range("F36").Value = "launch app"
$ws.Range("G36").Value = "https://adactinhotelapp.com/"
$ws.Range("G36").Style = "Hyperlink"

$ws.Range("E37").Value = "Step 3"
$ws.Range("F37").Value = "enter username in username textbox"
$ws.Range("G37").Value = "reyaz0806"
$ws.Range("H37").Value = "textbox should accept the text and display the same"

$ws.Range("E38").Value = "Step 4"
$ws.Range("F38").Value = "click login button"
$ws.Range("H38").Value = "Error should display below text `n`"Enter Password`""
$ws.Range("H38").WrapText = $true
$ws.Rows.Item(38).RowHeight = 46.8

# TC007 - Username filled, blank password
$ws.Range("A40").Value = "Adactin_Login_TS001_TC007"
$ws.Range("B40").Value = "To validdate login with blank username and password  field with data"
$ws.Range("B40").WrapText = $true
$ws.Range("C40").Value = "Adactin_Login_TS001"
$ws.Range("D40").Value = "user should have registered"
$ws.Range("D40").WrapText = $true
$ws.Range("E40").Value = "Step 1"
$ws.Range("F40").Value = "start  browser"
$ws.Range("G40").Value = "chrome"
$ws.Rows.Item(40).RowHeight = 93.6

$ws.Range("E41").Value = "Step 2"
$ws.Range("F41").Value = "launch app"
$ws.Range("G41").Value = "https://adactinhotelapp.com/"
$ws.Range("G41").Style = "Hyperlink"

$ws.Range("E42").Value = "Step 3"
$ws.Range("F42").Value = "enter password in password textbox"
$ws.Range("G42").Value = "reyaz123"
$ws.Range("H42").Value = "textbox should accept the text and display  as masked"

$ws.Range("E43").Value = "Step 4"
$ws.Range("F43").Value = "click login button"
$ws.Range("H43").Value = "Error should display below text `n`"Enter Username`""
$ws.Range("H43").WrapText = $true
$ws.Rows.Item(43).RowHeight = 46.8

# TC008 - Registered but unverified email
$ws.Range("A45").Value = "Adactin_Login_TS001_TC008"
$ws.Range("B45").Value = "To validdate login with registered username but not verified email"
$ws.Range("B45").WrapText = $true
$ws.Range("C45").Value = "Adactin_Login_TS001"
$ws.Range("D45").Value = "user should have registered and not verified his email"
$ws.Range("D45").WrapText = $true
$ws.Range("E45").Value = "Step 1"
$ws.Range("F45").Value = "start  browser"
$ws.Range("G45").Value = "chrome"
$ws.Range("H45").Value = "browser should get started"
$ws.Rows.Item(45).RowHeight = 93.6

$ws.Range("E46").Value = "Step 2"
$ws.Range("F46").Value = "launch app"
$ws.Range("G46").Value = "https://adactinhotelapp.com/"
$ws.Range("G46").Style = "Hyperlink"
$ws.Range("H46").Value = "app should get launched and land in home page"

$ws.Range("E47").Value = "Step 3"
$ws.Range("F47").Value = "enter username in username textbox"
$ws.Range("G47").Value = "murthi0602"
$ws.Range("H47").Value = "textbox should accept the text and display the same"

$ws.Range("E48").Value = "Step 4"
$ws.Range("F48").Value = "enter password in password textbox"
$ws.Range("G48").Value = "murthi123"
$ws.Range("H48").Value = "textbox should accept the text and display  as masked"

$ws.Range("E49").Value = "Step 5"
$ws.Range("F49").Value = "click login button"
$ws.Range("H49").Value = "login button should be clickable "

$ws.Range("H50").Value = "Error should display below text `n`"Error: Pending Email Verification`""
$ws.Range("H50").WrapText = $true
$ws.Rows.Item(50).RowHeight = 46.8

# Restore the cursor/selection position recorded in the edit.
$ws.Range("H54").Select() | Out-Null
